# ---------------------------------------------------------------------------
# Edit script: reproduces the authoral changes to
# docs/Apresentação_Comprod.pptx (commit "Add files via upload").
#
# Summary of changes applied:
#   Slide 1 (title slide):
#     - resize/move the big background rounded-rect behind the title
#     - resize/move the title placeholder
#     - drop the small cropped screenshot picture ("Imagem 2")
#     - move/enlarge the logo picture ("Imagem 4") into the freed-up space
#   Slide 5 ("Manual de instalação"):
#     - resize/move the two background rounded-rects
#     - widen the rounded-rect's corner radius (adjustment value)
#     - move/resize + uppercase the title text
#     - add a new subtitle placeholder with descriptive text
#   Slide 7 ("Help Desk"):
#     - same restyle as slide 5, with its own title/subtitle copy
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ===========================================================================
# Slide 1 - title slide
# ===========================================================================
$s1 = $p.Slides.Item(1)

# Big rounded rectangle behind the title block
$rr1 = $s1.Shapes.Item("Google Shape;292;p36")
$rr1.Left   = 100.88173298346457
$rr1.Top    = 137.49291238582677
$rr1.Width  = 528.4252014503937
$rr1.Height = 147.59449008897639

# Title placeholder ("COMPANHIA DE PRODUTIVIDADE / COMPROD")
$title1 = $s1.Shapes.Item("Google Shape;294;p36")
$title1.Left   = 85.39236070472441
$title1.Top    = 168.32173928346458
$title1.Width  = 550.9133606267717
$title1.Height = 86.19708641417323

# Drop the small cropped picture ("Imagem 2")
$img2 = $s1.Shapes.Item("Imagem 2")
$img2.Delete()

# Move / enlarge the logo picture ("Imagem 4") into the space freed up
$img4 = $s1.Shapes.Item("Imagem 4")
$img4.Left   = 456.4667816535433
$img4.Top    = 0.00007874015748031496
$img4.Width  = 256.42213448425196
$img4.Height = 86.2744102488189

# ===========================================================================
# Helper values shared by slide 5 / slide 7 restyle
# ===========================================================================

function Restyle-InstallHelpSlide($slide, $titleText, $subtitleText) {
    # First (right-hand) rounded rectangle -> becomes the full-width card
    # behind the title.
    $rr1 = $slide.Shapes.Item(1)
    $rr1.Left   = 95.78740317480316
    $rr1.Top    = 72.86378102755906
    $rr1.Width  = 528.4252014503937
    $rr1.Height = 197.88188936377952

    # Second (left-hand, tall) rounded rectangle -> becomes a thin footer
    # band, also gets a bigger corner radius.
    $rr2 = $slide.Shapes.Item(2)
    $rr2.Left   = 95.78740317480316
    $rr2.Top    = 284.13621527244095
    $rr2.Width  = 528.4252014503937
    $rr2.Height = 47.9999981
    $rr2.Adjustments.Item(1) = 0.2028

    # Title placeholder: reflow + uppercase text
    $title = $slide.Shapes.Item(3)
    $title.Left   = 119.55543137086615
    $title.Top    = 140.76535803070868
    $title.Width  = 480.88890083779523
    $title.Height = 62.078741157480316
    $title.TextFrame.TextRange.Text = $titleText

    # New subtitle placeholder, inherited from the slide's layout
    $slide.CustomLayout = $slide.CustomLayout
    $subtitle = $slide.Shapes.Item(4)
    $subtitle.Name   = "Google Shape;376;p41"
    $subtitle.Left   = 95.78740317480316
    $subtitle.Top    = 291.0220490440945
    $subtitle.Width  = 528.4252014503937
    $subtitle.Height = 34.22834595669291
    $subtitle.TextFrame.TextRange.Text = $subtitleText
}

# ===========================================================================
# Slide 5 - "Manual de instalação"
# ===========================================================================
$s5 = $p.Slides.Item(5)
Restyle-InstallHelpSlide $s5 "MANUAL DE INSTALAÇÃO" "Passo a passo para instalar e rodar o arduino no seu ambiente"

# ===========================================================================
# Slide 7 - "Help Desk"
# ===========================================================================
$s7 = $p.Slides.Item(7)
Restyle-InstallHelpSlide $s7 "HELPDESK " "Passo a Passo de como funciona e qual Ferramenta que escolhemos"

# Slide 7's subtitle sits slightly further right than slide 5's.
$subtitle7 = $s7.Shapes.Item("Google Shape;376;p41")
$subtitle7.Left = 96.6363792527559

Write-Output "edit complete"
